# "Generate Report for Archive"
#
# The localization-status report is regenerated: the file
# 5531544a-1ac6-4ce0-8c4b-a8862c0be44b.md has moved back into translation
# (status flips from "Ready for handoff" back to "In Translation"), which
# re-sorts it to the top of each report table; the other three files
# (b05fe1ca, c627f309, 52abe100) keep their own data untouched but shift
# down one row to make room.

$wb = $excel.ActiveWorkbook

# Final row order (row 2..5) for every sheet after the regenerate.
$order = @(
  '5531544a-1ac6-4ce0-8c4b-a8862c0be44b',
  'b05fe1ca-9a0e-4e45-9c36-95499bb03ada',
  'c627f309-323a-4055-9f70-c1e46d0d99dd',
  '52abe100-60f7-4f8b-93a5-08fb336a0d6e'
)

# Github blob URLs behind each hyperlink - keyed by file guid, independent of row.
$urls = @{
  'b05fe1ca-9a0e-4e45-9c36-95499bb03ada' = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2a5931f408b0281fec078443c0f9e0a7013ad78/e2e/b05fe1ca-9a0e-4e45-9c36-95499bb03ada.md'
  'c627f309-323a-4055-9f70-c1e46d0d99dd' = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a2a5931f408b0281fec078443c0f9e0a7013ad78/e2e/c627f309-323a-4055-9f70-c1e46d0d99dd.md'
  '52abe100-60f7-4f8b-93a5-08fb336a0d6e' = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1a5a6dce2aa6601ef5d9199418c7521ea6ab711e/e2e/52abe100-60f7-4f8b-93a5-08fb336a0d6e.md'
  '5531544a-1ac6-4ce0-8c4b-a8862c0be44b' = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/029c3cd784cc3571244e6461654ceee54ac4e02b/e2e/5531544a-1ac6-4ce0-8c4b-a8862c0be44b.md'
}

# ---------------------------------------------------------------------------
# Sheet "Overview" - columns: A File Name, B Path And Name, C Extension,
# D Publish URL, E zh-cn, F de-de, G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$overviewData = @{
  'b05fe1ca-9a0e-4e45-9c36-95499bb03ada' = @{ G = '2017-02-21 02:23:25'; Status = 'In Translation' }
  'c627f309-323a-4055-9f70-c1e46d0d99dd' = @{ G = '2017-02-21 02:23:25'; Status = 'In Translation' }
  '52abe100-60f7-4f8b-93a5-08fb336a0d6e' = @{ G = '2017-02-21 02:25:55'; Status = 'Ready for handoff' }
  '5531544a-1ac6-4ce0-8c4b-a8862c0be44b' = @{ G = '2017-02-21 02:25:06'; Status = 'In Translation' }
}

$wsOverview = $wb.Worksheets.Item('Overview')
$wsOverview.Hyperlinks.Delete()

for ($i = 0; $i -lt $order.Count; $i++) {
  $guid = $order[$i]
  $row = $i + 2
  $info = $overviewData[$guid]

  $wsOverview.Cells.Item($row, 1).Value = "$guid.md"
  $wsOverview.Cells.Item($row, 3).Value = '.md'
  $wsOverview.Cells.Item($row, 5).Value = $info.Status
  $wsOverview.Cells.Item($row, 6).Value = $info.Status
  $wsOverview.Cells.Item($row, 7).Value = $info.G

  $wsOverview.Hyperlinks.Add(
    $wsOverview.Cells.Item($row, 2),
    $urls[$guid],
    "",
    "",
    "e2e\$guid.md"
  )
}

# ---------------------------------------------------------------------------
# Sheet "zh-cn" / "de-de" - columns: A Source File Name, B File Extension,
# C Status, D Source Path, E Priority, F Content Duplicate,
# G Latest Handoff File, H Latest Handoff Datetime, ...
# ---------------------------------------------------------------------------
$langData = @{
  'zh-cn' = @{
    'b05fe1ca-9a0e-4e45-9c36-95499bb03ada' = @{ G = 'b05fe1ca-9a0e-4e45-9c36-95499bb03ada.9dcb0024235953b2a43f551bdeb92bba837c8600.zh-cn.xlf'; H = '2017-02-21 02:23:08'; Status = 'In Translation' }
    'c627f309-323a-4055-9f70-c1e46d0d99dd' = @{ G = 'c627f309-323a-4055-9f70-c1e46d0d99dd.dacdfa0e42d497832252e217de8e75f5256da4e8.zh-cn.xlf'; H = '2017-02-21 02:23:08'; Status = 'In Translation' }
    '52abe100-60f7-4f8b-93a5-08fb336a0d6e' = @{ G = '52abe100-60f7-4f8b-93a5-08fb336a0d6e.d1cbf8b7f020f04b71352495d83a62ed707f7a35.zh-cn.xlf'; H = '2017-02-21 02:25:40'; Status = 'Ready for handoff' }
    '5531544a-1ac6-4ce0-8c4b-a8862c0be44b' = @{ G = '5531544a-1ac6-4ce0-8c4b-a8862c0be44b.94b5a1cd267f1f0ef6706966dce17c5acd4127e1.zh-cn.xlf'; H = '2017-02-21 02:24:49'; Status = 'In Translation' }
  }
  'de-de' = @{
    'b05fe1ca-9a0e-4e45-9c36-95499bb03ada' = @{ G = 'b05fe1ca-9a0e-4e45-9c36-95499bb03ada.9dcb0024235953b2a43f551bdeb92bba837c8600.de-de.xlf'; H = '2017-02-21 02:23:25'; Status = 'In Translation' }
    'c627f309-323a-4055-9f70-c1e46d0d99dd' = @{ G = 'c627f309-323a-4055-9f70-c1e46d0d99dd.dacdfa0e42d497832252e217de8e75f5256da4e8.de-de.xlf'; H = '2017-02-21 02:23:25'; Status = 'In Translation' }
    '52abe100-60f7-4f8b-93a5-08fb336a0d6e' = @{ G = '52abe100-60f7-4f8b-93a5-08fb336a0d6e.d1cbf8b7f020f04b71352495d83a62ed707f7a35.de-de.xlf'; H = '2017-02-21 02:25:55'; Status = 'Ready for handoff' }
    '5531544a-1ac6-4ce0-8c4b-a8862c0be44b' = @{ G = '5531544a-1ac6-4ce0-8c4b-a8862c0be44b.94b5a1cd267f1f0ef6706966dce17c5acd4127e1.de-de.xlf'; H = '2017-02-21 02:25:06'; Status = 'In Translation' }
  }
}

foreach ($langName in @('zh-cn', 'de-de')) {
  $ws = $wb.Worksheets.Item($langName)
  $ws.Hyperlinks.Delete()
  $data = $langData[$langName]

  for ($i = 0; $i -lt $order.Count; $i++) {
    $guid = $order[$i]
    $row = $i + 2
    $info = $data[$guid]

    $ws.Cells.Item($row, 1).Value = "$guid.md"
    $ws.Cells.Item($row, 3).Value = $info.Status
    $ws.Cells.Item($row, 7).Value = $info.G
    $ws.Cells.Item($row, 8).Value = $info.H

    $ws.Hyperlinks.Add(
      $ws.Cells.Item($row, 1),
      $urls[$guid],
      "",
      "",
      "$guid.md"
    )
  }
}
